# Horarios actualizados Linea 141 - 481
# Updates the three schedule sheets with the newly scraped data.

$wb = $excel.ActiveWorkbook

$nuevaHora = "04:21:55"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 1).Value2 = "Última actualización: $nuevaHora"
$ws1.Cells.Item(3, 1).Value2 = "Total filas: 10"

$sheet1Data = @(
    ,@($nuevaHora, "04:46", "215_EL PELIGRO", 25,  "LP1912")
    ,@($nuevaHora, "04:53", "11_ETCHEVERRY",  32,  "LP1912")
    ,@($nuevaHora, "05:11", "17_ROMERO",       50,  "LP1912")
    ,@($nuevaHora, "05:21", "23_HERNANDEZ",    60,  "LP1912")
    ,@($nuevaHora, "05:31", "81_EL PELIGRO",   70,  "LP1912")
    ,@($nuevaHora, "05:38", "14_ABASTO",       77,  "LP1912")
    ,@($nuevaHora, "05:51", "17_ROMERO",       90,  "LP1912")
    ,@($nuevaHora, "06:00", "16_SANTA ANA",    99,  "LP1912")
    ,@($nuevaHora, "06:03", "10_OLMOS",        102, "LP1912")
    ,@($nuevaHora, "06:10", "215A_EL PATO",    109, "LP1912")
)

$row = 6
foreach ($r in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value2 = $r[0]
    $ws1.Cells.Item($row, 2).Value2 = $r[1]
    $ws1.Cells.Item($row, 3).Value2 = $r[2]
    $ws1.Cells.Item($row, 4).Value2 = $r[3]
    $ws1.Cells.Item($row, 5).Value2 = $r[4]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 1).Value2 = "Última actualización: $nuevaHora"
$ws2.Cells.Item(3, 1).Value2 = "Total filas: 2"

$sheet2Data = @(
    ,@($nuevaHora, "04:46", "215_EL PELIGRO", 25,  "LP1912")
    ,@($nuevaHora, "06:10", "215A_EL PATO",    109, "LP1912")
)

$row = 6
foreach ($r in $sheet2Data) {
    $ws2.Cells.Item($row, 1).Value2 = $r[0]
    $ws2.Cells.Item($row, 2).Value2 = $r[1]
    $ws2.Cells.Item($row, 3).Value2 = $r[2]
    $ws2.Cells.Item($row, 4).Value2 = $r[3]
    $ws2.Cells.Item($row, 5).Value2 = $r[4]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 (only the timestamp refreshes, no data rows)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2, 1).Value2 = "Última actualización: $nuevaHora"
